$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$r1 = $t.Rows.Item(1)
Write-Host "HeadingFormat before:" $r1.HeadingFormat
$r1.HeadingFormat = 0
Write-Host "HeadingFormat after:" $r1.HeadingFormat
